$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Col8a1-Itga1 LR-pair results (new Sending/Target cluster cross: ECs, FAPs, sCs)
$data = @(
    ,("ECs", "Col8a1", "Itga1", "ECs", 2, 0.6666666666666666, 3.346688666666667, 10.040066, 0.1019109290125008, 0.1019109290125008, 3, 1, 22.11318066666666, 66.339542, 0.1993293533530854, 0.1993293533530854, 74.00593112108577, 666.053380089772, 0.02031383957967397, 0.02031383957967397)
    ,("ECs", "Col8a1", "Itga1", "FAPs", 2, 0.6666666666666666, 3.346688666666667, 10.040066, 0.1019109290125008, 0.1019109290125008, 3, 1, 7.622575, 22.867725, 0.06871028498970018, 0.06871028498970018, 25.51038536331667, 229.59346826985, 0.007002328976014031, 0.007002328976014031)
    ,("ECs", "Col8a1", "Itga1", "sCs", 2, 0.6666666666666666, 3.346688666666667, 10.040066, 0.1019109290125008, 0.1019109290125008, 3, 1, 81.20214833333334, 243.606445, 0.7319603616572145, 0.7319603616572145, 271.7583095361522, 2445.82478582537, 0.07459476045681276, 0.07459476045681276)
    ,("FAPs", "Col8a1", "Itga1", "ECs", 3, 1, 26.21646566666666, 78.649397, 0.7983247435368442, 0.7983247435368442, 3, 1, 22.11318066666666, 66.339542, 0.1993293533530854, 0.1993293533530854, 579.7294417284637, 5217.564975556174, 0.1591295548949669, 0.1591295548949669)
    ,("FAPs", "Col8a1", "Itga1", "FAPs", 3, 1, 26.21646566666666, 78.649397, 0.7983247435368442, 0.7983247435368442, 3, 1, 7.622575, 22.867725, 0.06871028498970018, 0.06871028498970018, 199.8369757790917, 1798.532782011825, 0.05485312064274588, 0.05485312064274588)
    ,("FAPs", "Col8a1", "Itga1", "sCs", 3, 1, 26.21646566666666, 78.649397, 0.7983247435368442, 0.7983247435368442, 3, 1, 81.20214833333334, 243.606445, 0.7319603616572145, 0.7319603616572145, 2128.833333840407, 19159.50000456366, 0.5843420679991315, 0.5843420679991315)
    ,("sCs", "Col8a1", "Itga1", "ECs", 3, 1, 3.276195666666666, 9.828586999999999, 0.09976432745065499, 0.09976432745065499, 3, 1, 22.11318066666666, 66.339542, 0.1993293533530854, 0.1993293533530854, 72.44710667635043, 652.0239600871539, 0.01988595887844452, 0.01988595887844452)
    ,("sCs", "Col8a1", "Itga1", "FAPs", 3, 1, 3.276195666666666, 9.828586999999999, 0.09976432745065499, 0.09976432745065499, 3, 1, 7.622575, 22.867725, 0.06871028498970018, 0.06871028498970018, 24.97304718384166, 224.757424654575, 0.006854835370940273, 0.006854835370940273)
    ,("sCs", "Col8a1", "Itga1", "sCs", 3, 1, 3.276195666666666, 9.828586999999999, 0.09976432745065499, 0.09976432745065499, 3, 1, 81.20214833333334, 243.606445, 0.7319603616572145, 0.7319603616572145, 266.0341264936906, 2394.307138443215, 0.0730235332012702, 0.0730235332012702)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
